$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '49.505.33'
$ws.Range("E2").Value = '  -0.81%  '

$ws.Range("D3").Value = '2.637.85'
$ws.Range("E3").Value = '  -0.44%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.60'
$ws.Range("E5").Value = '  -1.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '325.95'
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("E7").Value = '  -0.99%  '

$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("E9").Value = '  -0.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.64'
$ws.Range("E10").Value = '  -3.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.00'
$ws.Range("E11").Value = '  -0.87%  '

$ws.Range("E12").Value = '  -0.59%  '

$ws.Range("E13").Value = '  +1.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.60'
$ws.Range("E14").Value = '  +3.53%  '

$ws.Range("D15").Value = '3.051.32'
$ws.Range("E15").Value = '  -0.33%  '

$ws.Range("D16").Value = '2.629.71'
$ws.Range("E16").Value = '  -0.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.857'
$ws.Range("E17").Value = '  -1.78%  '

$ws.Range("D18").Value = '49.493.65'
$ws.Range("E18").Value = '  -0.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.38'
$ws.Range("E19").Value = '  +1.93%  '

$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.67'
$ws.Range("E21").Value = '  -1.79%  '

$ws.Range("D22").Value = '0.0₃0949'
$ws.Range("E22").Value = '  -0.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '268.27'
$ws.Range("E23").Value = '  -3.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.09'
$ws.Range("E24").Value = '  -4.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.56'
$ws.Range("E25").Value = '  -0.98%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.02'
$ws.Range("E27").Value = '  -2.95%  '

$ws.Range("E28").Value = '  +1.44%  '

$ws.Range("E29").Value = '  -1.12%  '

$ws.Range("E30").Value = '  -2.62%  '

$ws.Range("E31").Value = '  -3.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.59'
$ws.Range("E32").Value = '  -1.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.48'
$ws.Range("E33").Value = '  +0.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0817'
$ws.Range("E34").Value = '  +0.75%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.12'
$ws.Range("E35").Value = '  -1.67%  '

$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.95'
$ws.Range("E37").Value = '  +2.08%  '

$ws.Range("E38").Value = '  -2.28%  '

$ws.Range("E39").Value = '  -0.23%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '129.53'
$ws.Range("E40").Value = '  +4.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.31'
$ws.Range("E41").Value = '  +3.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.01'
$ws.Range("E42").Value = '  +4.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0341'
$ws.Range("E43").Value = '  +8.22%  '

$ws.Range("E44").Value = '  -1.02%  '

$ws.Range("D45").Value = '2.062.17'
$ws.Range("E45").Value = '  -0.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.30'
$ws.Range("E46").Value = '  -0.75%  '

$ws.Range("E47").Value = '  +7.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.16'
$ws.Range("E48").Value = '  -7.26%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.89'
$ws.Range("E49").Value = '  -2.57%  '

$ws.Range("E50").Value = '  -3.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.67'
$ws.Range("E51").Value = '  -0.91%  '
